$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - shift header labels left by one column and add TotalCasesPer1M
$ws.Range("A1").Value = "TotalCases"
$ws.Range("B1").Value = "NewCases"
$ws.Range("C1").Value = "TotalDeaths"
$ws.Range("D1").Value = "NewDeaths"
$ws.Range("E1").Value = "TotalRecovered"
$ws.Range("F1").Value = "NewRecovered"
$ws.Range("G1").Value = "ActiveCases"
$ws.Range("H1").Value = "Critical"
$ws.Range("I1").Value = "TotalCasesPer1M"

# Data rows: drop the Continent column, shift remaining data left, refresh with latest figures
# Row 2
$ws.Range("A2").Value = 131591034
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 1692184
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = 127109905
$ws.Range("F2").Value = 665
$ws.Range("G2").Value = 2788945
$ws.Range("H2").Value = 6709
$ws.Range("I2").ClearContents()
# Row 3
$ws.Range("A3").Value = 221465849
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 1553248
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = 205656562
$ws.Range("F3").Value = 308
$ws.Range("G3").Value = 14256039
$ws.Range("H3").Value = 14733
$ws.Range("I3").ClearContents()
# Row 4
$ws.Range("A4").Value = 253230160
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = 2100520
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = 248563147
$ws.Range("F4").Value = 5389
$ws.Range("G4").Value = 2566493
$ws.Range("H4").Value = 4520
$ws.Range("I4").ClearContents()
# Row 5
$ws.Range("A5").Value = 69859756
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = 1365807
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = 66650237
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = 1843712
$ws.Range("H5").Value = 8953
$ws.Range("I5").ClearContents()
# Row 6
$ws.Range("A6").Value = 14843156
$ws.Range("B6").Value = 806
$ws.Range("C6").Value = 32737
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = 14578860
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = 231559
$ws.Range("H6").Value = 40
$ws.Range("I6").ClearContents()
# Row 7
$ws.Range("A7").Value = 12860287
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = 258884
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = 12090790
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = 510613
$ws.Range("H7").Value = 529
$ws.Range("I7").ClearContents()
# Row 8
$ws.Range("A8").Value = 721
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = 15
$ws.Range("D8").ClearContents()
$ws.Range("E8").Value = 706
$ws.Range("F8").ClearContents()
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").ClearContents()
